$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new I1:J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns: I0 and IF values per row
$data = @(
    @(7, 7),
    @(6, 7),
    @(7, 8),
    @(3, 3),
    @(5, 6),
    @(4, 5),
    @(5, 6),
    @(10, 10),
    @(10, 10),
    @(7, 8),
    @(6, 7),
    @(3, 5),
    @(5, 7),
    @(5, 7),
    @(9, 9),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
